$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 269.7619
$ws.Range("I33").Value = 408.72726
$ws.Range("K33").Value = 408.72726
$ws.Range("M33").Value = -179.72726
$ws.Range("H70").Value = 1477.1578
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1477.1578
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4431.4734
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -4971.4734
$ws.Range("H73").Value = 1477.1578
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1477.1578
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4431.4734
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -6303.4734
$ws.Range("H86").Value = 412044.34
$ws.Range("H89").Value = 412044.34
$ws.Range("H112").Value = 2500
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 2780
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 8340
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -10556
$ws.Range("H118").Value = 676
$ws.Range("I118").Value = 676
$ws.Range("K118").Value = 2028
$ws.Range("M118").Value = -371
$ws.Range("H132").Value = 905.1
$ws.Range("I132").Value = 842.7174
$ws.Range("K132").Value = 2528.1522
$ws.Range("M132").Value = 1.847800000000007
$ws.Range("H137").Value = 2303.5
$ws.Range("I137").Value = 1405.1666
$ws.Range("J137").Value = 2602.9443
$ws.Range("K137").Value = 4215.4998
$ws.Range("L137").Value = 7808.8329
$ws.Range("M137").Value = -1665.4998
$ws.Range("N137").Value = -12908.8329
$ws.Range("H138").Value = 2795.5334
$ws.Range("J138").Value = 2774.0571
$ws.Range("L138").Value = 8322.1713
$ws.Range("N138").Value = -18602.1713
$ws.Range("H140").Value = 57114.375
$ws.Range("J140").Value = 57114.375
$ws.Range("L140").Value = 57114.375
$ws.Range("N140").Value = -67474.375

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5004900
$ws.Range("J6").Value = 9800
$ws.Range("L6").Value = 9800
$ws.Range("N6").Value = -10146
$ws.Range("H32").Value = 3985.238
$ws.Range("I32").Value = 3582.027
$ws.Range("K32").Value = 3582.027
$ws.Range("M32").Value = -3295.027
$ws.Range("H45").Value = 1625.1177
$ws.Range("I45").Value = 1381.6666
$ws.Range("K45").Value = 1381.6666
$ws.Range("M45").Value = -1004.6666
$ws.Range("H74").Value = 1516.0952
$ws.Range("I74").Value = 1358.5
$ws.Range("K74").Value = 1358.5
$ws.Range("M74").Value = -484.5
$ws.Range("H77").Value = 1516.0952
$ws.Range("I77").Value = 1358.5
$ws.Range("K77").Value = 6792.5
$ws.Range("M77").Value = -2424.5
$ws.Range("H122").Value = 1416.4445
$ws.Range("I122").Value = 1649.8334
$ws.Range("K122").Value = 4949.5002
$ws.Range("M122").Value = -2499.5002
$ws.Range("H132").Value = 1535.4468
$ws.Range("I132").Value = 951.4474
$ws.Range("J132").Value = 4001.2222
$ws.Range("K132").Value = 2854.3422
$ws.Range("L132").Value = 12003.6666
$ws.Range("M132").Value = -324.3422
$ws.Range("N132").Value = -17063.6666

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2208.4443
$ws.Range("I20").Value = 2152.913
$ws.Range("K20").Value = 2152.913
$ws.Range("M20").Value = -1905.913
$ws.Range("H36").Value = 57041
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 57041
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 57041
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -58109
$ws.Range("H134").Value = 5555.263
$ws.Range("I134").Value = 5900.0625
$ws.Range("K134").Value = 17700.1875
$ws.Range("M134").Value = -15165.1875

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2411.8823
$ws.Range("I31").Value = 3068.1667
$ws.Range("K31").Value = 3068.1667
$ws.Range("M31").Value = -2773.1667
$ws.Range("H34").Value = 2411.8823
$ws.Range("I34").Value = 3068.1667
$ws.Range("K34").Value = 3068.1667
$ws.Range("M34").Value = -2866.1667
$ws.Range("H62").Value = 2079.6
$ws.Range("J62").Value = 1599.5
$ws.Range("L62").Value = 1599.5
$ws.Range("N62").Value = -2847.5
$ws.Range("H65").Value = 2079.6
$ws.Range("J65").Value = 1599.5
$ws.Range("L65").Value = 7997.5
$ws.Range("N65").Value = -14237.5
$ws.Range("H132").Value = 1866.4839
$ws.Range("I132").Value = 1011.2
$ws.Range("K132").Value = 3033.6
$ws.Range("M132").Value = -503.6000000000004
$ws.Range("H134").Value = 2002.2646
$ws.Range("I134").Value = 1769.2667
$ws.Range("K134").Value = 5307.800099999999
$ws.Range("M134").Value = -2772.800099999999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 4976.4
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4976.4
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 14929.2
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -15267.2
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H107").Value = 843.82355
$ws.Range("I107").Value = 612.5714
$ws.Range("J107").Value = 1005.7
$ws.Range("K107").Value = 1837.7142
$ws.Range("L107").Value = 3017.1
$ws.Range("M107").Value = 82.28579999999988
$ws.Range("N107").Value = -6857.1
$ws.Range("H115").Value = 4085.0715
$ws.Range("J115").Value = 5986.625
$ws.Range("L115").Value = 17959.875
$ws.Range("N115").Value = -20309.875
$ws.Range("H122").Value = 879.5789
$ws.Range("J122").Value = 1127.9166
$ws.Range("L122").Value = 10151.2494
$ws.Range("N122").Value = -15051.2494
$ws.Range("H131").Value = 5690615.5
$ws.Range("J131").Value = 9393.415000000001
$ws.Range("L131").Value = 28180.245
$ws.Range("N131").Value = -38260.245
$ws.Range("H134").Value = 2452.375
$ws.Range("I134").Value = 1250.6666
$ws.Range("J134").Value = 3997.4285
$ws.Range("K134").Value = 3751.9998
$ws.Range("L134").Value = 11992.2855
$ws.Range("M134").Value = 1318.0002
$ws.Range("N134").Value = -22132.2855

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7399.8
$ws.Range("I70").Value = 6999.5
$ws.Range("K70").Value = 6999.5
$ws.Range("M70").Value = -6729.5
$ws.Range("H73").Value = 7399.8
$ws.Range("I73").Value = 6999.5
$ws.Range("K73").Value = 6999.5
$ws.Range("M73").Value = -6063.5
$ws.Range("H93").Value = 29499.666
$ws.Range("J93").Value = 29499.666
$ws.Range("L93").Value = 29499.666
$ws.Range("N93").Value = -33243.666
$ws.Range("H107").Value = 1214.8334
$ws.Range("I107").Value = 64.666664
$ws.Range("K107").Value = 64.666664
$ws.Range("M107").Value = 1855.333336
$ws.Range("H132").Value = 1242995.4
$ws.Range("I132").Value = 1604191.9
$ws.Range("K132").Value = 4812575.699999999
$ws.Range("M132").Value = -4810045.699999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1663.7142
$ws.Range("I82").Value = 1500.5
$ws.Range("J82").Value = 1881.3334
$ws.Range("K82").Value = 1500.5
$ws.Range("L82").Value = 1881.3334
$ws.Range("M82").Value = -1139.5
$ws.Range("N82").Value = -2603.3334
$ws.Range("H85").Value = 1663.7142
$ws.Range("I85").Value = 1500.5
$ws.Range("J85").Value = 1881.3334
$ws.Range("K85").Value = 1500.5
$ws.Range("L85").Value = 1881.3334
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -4377.3334
$ws.Range("H132").Value = 2100.52
$ws.Range("I132").Value = 1250.8889
$ws.Range("J132").Value = 4285.2856
$ws.Range("K132").Value = 3752.6667
$ws.Range("L132").Value = 12855.8568
$ws.Range("M132").Value = -1222.6667
$ws.Range("N132").Value = -17915.8568
$ws.Range("H136").Value = 2498.7188
$ws.Range("I136").Value = 1447.12
$ws.Range("K136").Value = 4341.36
$ws.Range("M136").Value = -1791.36

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 119429
$ws.Range("J46").Value = 119429
$ws.Range("L46").Value = 119429
$ws.Range("N46").Value = -119891
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H132").Value = 2160.5557
$ws.Range("I132").Value = 1304.9375
$ws.Range("J132").Value = 3405.0908
$ws.Range("K132").Value = 3914.8125
$ws.Range("L132").Value = 10215.2724
$ws.Range("M132").Value = -1384.8125
$ws.Range("N132").Value = -15275.2724
$ws.Range("H134").Value = 119429
$ws.Range("J134").Value = 119429
$ws.Range("L134").Value = 358287
$ws.Range("N134").Value = -363357
$ws.Range("H136").Value = 20577816
$ws.Range("I136").Value = 26456194
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 79368582
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -79366032
$ws.Range("N136").Value = -15600
